$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark that currently sits after "a random" ---
try {
    $existing = $d.Bookmarks("_GoBack")
    $existing.Delete()
} catch {
    # no pre-existing bookmark with that name - nothing to remove
}

# --- 2. Remove ", e.g. ‘Pthreads’ and ‘OpenMP’" that follows "‘multithreading’" ---
$rng = $d.Content
$searchText = ", e.g. ‘Pthreads’ and ‘OpenMP’"
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = ""
    # --- 3. Re-create the "_GoBack" bookmark at the now-collapsed spot ---
    $d.Bookmarks.Add("_GoBack", $rng)
}

# --- 4. Insert a new "Thread starvation" list item right after the "Race condition" item ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Race condition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.End = $rng2.End + 1
    $racePara = $rng2.Paragraphs.Item(1)

    $raceIdx = 0
    $count = 0
    foreach ($p in $d.Paragraphs) {
        $count = $count + 1
        if ($p.Range.Start -eq $racePara.Range.Start) {
            $raceIdx = $count
        }
    }

    $insAt = $racePara.Range
    $insAt.Collapse(0)
    $insAt.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($raceIdx + 1)
    $newPara.Range.Text = "Thread starvation"
}
